$d = $word.ActiveDocument

# --- Paragraph 1: "NAME:" -> add the student's name ---------------------
$p1 = $d.Paragraphs(1).Range
$p1.InsertAfter(" Harsh Deep Keshari")

# --- Paragraph 2: "ENROLLMENT NO.:" -> add the enrollment number --------
$p2 = $d.Paragraphs(2).Range
$p2.InsertAfter(" A7605219086")

# --- Paragraph 3: "Description of  Class Case Study:" -> add the topic -
$p3 = $d.Paragraphs(3).Range
$p3.InsertAfter(" Darktable")

# --- New paragraph 4: the Darktable description -------------------------
$p3again = $d.Paragraphs(3).Range
$p3again.InsertParagraphAfter()

$p4 = $d.Paragraphs(4)
$p4.Range.Text = "Darktable is an open-source photography workflow application and raw developer. A virtual lighttable and darkroom for photographers. It manages your digital negatives in a database, lets you view them through a zoomable lighttable and enables you to develop raw images and enhance them."
$p4.Range.Font.Size = 18
$p4.Range.Font.SizeBi = 18
$p4.Format.Alignment = 3
